$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 789 (shifts existing rows 789-830 down to 790-831)
$ws.Rows.Item(789).Insert()

# The date value "2026/02/13" looks like a date, so force text entry the same
# way Excel requires (format as text, enter value, then restore default
# "Normal" style so the stored cell matches its neighbours with no explicit
# style index).
$ws.Cells.Item(789, 1).NumberFormat = "@"
$ws.Cells.Item(789, 1).Value = "2026/02/13"
$ws.Cells.Item(789, 1).Style = "Normal"

$ws.Cells.Item(789, 2).Value = "金"
$ws.Cells.Item(789, 3).Value = 20
$ws.Cells.Item(789, 4).Value = 22
